# Improve DOCX export for investigation tiac document export
#
# The "Transféré à" and "Date de réception" paragraphs were built out of
# several fragmented runs (one run per Jinja2 template token). Collapse
# each of those paragraphs back down to a single run holding the full
# template text, and extend the "transfered_to" placeholder with a
# `or '-'` fallback so empty values render as "-".

$d = $word.ActiveDocument

function Set-ParagraphText($paragraphIndex, $newText) {
    $p = $d.Paragraphs($paragraphIndex)
    $start = $p.Range.Start
    $end = $p.Range.End - 1   # exclude the trailing paragraph mark
    $r = $d.Range($start, $end)
    $r.Text = $newText
}

# Locate the two target paragraphs by their current (fragmented) text so
# the script stays correct even if paragraph numbering shifts.
$transfereIndex = $null
$dateReceptionIndex = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Transf*r* *: {{ object.transfered_to }}*") {
        $transfereIndex = $i
    }
    if ($t -like "Date de r*ception*la DD(ETS)PP*date_reception.strftime*") {
        $dateReceptionIndex = $i
    }
}

if ($null -eq $transfereIndex) {
    throw "Could not locate the 'Transféré à' paragraph"
}
if ($null -eq $dateReceptionIndex) {
    throw "Could not locate the 'Date de réception' paragraph"
}

Set-ParagraphText $transfereIndex "Transféré à : {{ object.transfered_to  or '-'  }}"
Set-ParagraphText $dateReceptionIndex 'Date de réception à la DD(ETS)PP : {{object.date_reception.strftime("%Y-%m-%d %H:%M") }}'

Write-Output "Updated paragraphs $transfereIndex and $dateReceptionIndex"
